$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.726.98'
$ws.Range('E2').Value = '  +0.28%  '
$ws.Range('D3').Value = '2.458.75'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '573.58'
$ws.Range('E5').Value = '  -0.63%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.29'
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '2.457.59'
$ws.Range('E9').Value = '  +0.81%  '
$ws.Range('E10').Value = '  +1.12%  '
$ws.Range('E11').Value = '  +1.64%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').Value = '  +0.44%  '
$ws.Range('E13').Value = '  +1.03%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '28.80'
$ws.Range('E14').Value = '  +1.52%  '
$ws.Range('E15').Value = '  -0.62%  '
$ws.Range('D16').Value = '2.905.40'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('D17').Value = '62.780.44'
$ws.Range('E17').Value = '  +0.52%  '
$ws.Range('D18').Value = '2.461.57'
$ws.Range('E18').Value = '  +1.02%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '8.00'
$ws.Range('E19').Value = '  +2.65%  '
$ws.Range('E20').Value = '  +0.85%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '326.23'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('B22').Value = 'SuiNetwork'
$ws.Range('C22').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '2.22'
$ws.Range('E22').Value = '  +10.34%  '
$ws.Range('B23').Value = 'Polkadot'
$ws.Range('C23').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.13'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '10.05'
$ws.Range('E25').Value = '  +18.37%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '65.57'
$ws.Range('E26').Value = '  +0.49%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '651.60'
$ws.Range('E27').Value = '  +2.79%  '
$ws.Range('D29').Value = '0.0₃0976'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('E30').Value = '  -12.18%  '
$ws.Range('E31').Value = '  +3.21%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.97'
$ws.Range('E32').Value = '  -2.49%  '
$ws.Range('E33').Value = '  -0.84%  '
$ws.Range('E34').Value = '  -3.83%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.998'
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +2.20%  '
$ws.Range('E37').Value = '  +0.35%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '152.53'
$ws.Range('E39').Value = '  -0.76%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '18.66'
$ws.Range('E40').Value = '  +0.61%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.35'
$ws.Range('E41').Value = '  -1.63%  '
$ws.Range('E42').Value = '  +1.18%  '
$ws.Range('E43').Value = '  -1.80%  '
$ws.Range('D44').Value = '0.0₆0313'
$ws.Range('E44').Value = '  -69.90%  '
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '153.05'
$ws.Range('E46').Value = '  +6.27%  '
$ws.Range('E47').Value = '  +1.40%  '
$ws.Range('E48').Value = '  -0.19%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.605'
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('B50').Value = 'InjectiveProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '20.23'
$ws.Range('E50').Value = '  -0.99%  '
$ws.Range('E51').Value = '  -0.17%  '
